# Generate Report for Handback
# Updates the handback-status report with freshly generated timestamps
# for the cc58bce8-649e-473c-9bf0-bd9dfda0c66b.md file row.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
#     cc58bce8-649e-473c-9bf0-bd9dfda0c66b.md row (row 6) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G6").Value = "2017-01-03 07:36:37"

# --- zh-cn sheet: "Correspond Handoff Datetime" (H) and
#     "Correspond Handback DateTime" (L) for row 6 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H6").Value = "2017-01-03 07:36:25"
$wsZhCn.Range("L6").Value = "2017-01-03 07:36:57"

# --- de-de sheet: "Correspond Handoff Datetime" (H) and
#     "Correspond Handback DateTime" (L) for row 6 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H6").Value = "2017-01-03 07:36:37"
$wsDeDe.Range("L6").Value = "2017-01-03 07:37:09"
